$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (correct-answer mark value)
$ws.Range("B11").Value = 5

# Update "Total" row (total correct marks and Corr/total marks text)
$ws.Range("B12").Value = 70
$ws.Range("E12").Value = "70/140"
